# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" sheets:
#   F2: 257 -> 258
#   F4: 159 -> 161

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 258
    $ws.Range("F4").Value = 161
}
